# This script rewrites the counters summary sheet (rows 18-44) to add
# latitude/longitude format-consistency rows, a corrected type_id /
# type_name block, and new product/service enumeration error rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('latitude', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('latitude', 'Format Consistency', 'Error: Value '''' does not match regex ''^(\+|-)?(?:90(?:(?:\.0{1,6})?)|(?:[0-9]|[1-8][0-9])(?:(?:\.[0-9]{1,6})?))$'''),
    @('latitude', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('latitude', 'Format Consistency', 'Error: Value '''' does not match regex ''^(\+|-)?(?:90(?:(?:\.0{1,6})?)|(?:[0-9]|[1-8][0-9])(?:(?:\.[0-9]{1,6})?))$'''),
    @('latitude', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('latitude', 'Format Consistency', 'Error: Value '''' does not match regex ''^(\+|-)?(?:90(?:(?:\.0{1,6})?)|(?:[0-9]|[1-8][0-9])(?:(?:\.[0-9]{1,6})?))$'''),
    @('longitude', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('longitude', 'Format Consistency', 'Error: Value '''' does not match regex ''^(\+|-)?(?:180(?:(?:\.0{1,6})?)|(?:[0-9]|[1-9][0-9]|1[0-7][0-9])(?:(?:\.[0-9]{1,6})?))$'''),
    @('longitude', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('longitude', 'Format Consistency', 'Error: Value '''' does not match regex ''^(\+|-)?(?:180(?:(?:\.0{1,6})?)|(?:[0-9]|[1-9][0-9]|1[0-7][0-9])(?:(?:\.[0-9]{1,6})?))$'''),
    @('longitude', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('longitude', 'Format Consistency', 'Error: Value '''' does not match regex ''^(\+|-)?(?:180(?:(?:\.0{1,6})?)|(?:[0-9]|[1-9][0-9]|1[0-7][0-9])(?:(?:\.[0-9]{1,6})?))$'''),
    @('type_id', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('type_id', 'Meta Compliance (data type)', 'Error: Value '''' is not an int. An int was expected'),
    @('type_id', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('type_id', 'Meta Compliance (data type)', 'Error: Value '''' is not an int. An int was expected'),
    @('type_id', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('type_id', 'Meta Compliance (data type)', 'Error: Value '''' is not an int. An int was expected'),
    @('type_name', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('type_name', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('type_name', 'Completeness of Mandatory fields', 'Error: Mandatory field is BLANK or NULL. A value is required.'),
    @('product1', 'Meta Compliance (enumeration)', 'Error: Value ''225'' is outside the enumeration set ''[''go card'',''Gold Coast go explore card'',''Visitor Information Pack'',''SEEQ card'']'''),
    @('product2', 'Meta Compliance (enumeration)', 'Error: Value ''773'' is outside the enumeration set ''[''go card'',''Gold Coast go explore card'',''Visitor Information Pack'',''SEEQ card'']'''),
    @('product3', 'Meta Compliance (enumeration)', 'Error: Value ''859'' is outside the enumeration set ''[''go card'',''Gold Coast go explore card'',''Visitor Information Pack'',''SEEQ card'']'''),
    @('service1', 'Meta Compliance (enumeration)', 'Error: Value ''87'' is outside the enumeration set ''[''Buy'', ''Change expiry date'', ''Top up'', ''Refund'', ''Register'']'''),
    @('service2', 'Meta Compliance (enumeration)', 'Error: Value ''495'' is outside the enumeration set ''[''Buy'', ''Change expiry date'', ''Top up'', ''Refund'', ''Register'']'''),
    @('service3', 'Meta Compliance (enumeration)', 'Error: Value ''546'' is outside the enumeration set ''[''Buy'', ''Change expiry date'', ''Top up'', ''Refund'', ''Register'']''')
)

$startRow = 18
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

